$wb = $excel.ActiveWorkbook

# row 32 on ALC (-2206,25 +2206,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1600.1666
$ws.Range("J32").Value = 1320.2
$ws.Range("L32").Value = 1320.2
$ws.Range("N32").Value = -1972.2

# row 51 on ALC (-3155,25 +3155,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3979.8696
$ws.Range("I51").Value = 4466.6665
$ws.Range("J51").Value = 3067.125
$ws.Range("K51").Value = 4466.6665
$ws.Range("L51").Value = 3067.125
$ws.Range("M51").Value = -3982.6665
$ws.Range("N51").Value = -4035.125

# row 58 on ALC (-3513,23 +3513,26)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1757.625
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9300

# row 62 on ALC (-3703,25 +3706,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6302.8
$ws.Range("I62").Value = 5737.875
$ws.Range("J62").Value = 6948.4287
$ws.Range("K62").Value = 5737.875
$ws.Range("L62").Value = 6948.4287
$ws.Range("M62").Value = -5113.875
$ws.Range("N62").Value = -8196.4287

# row 65 on ALC (-3856,25 +3859,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 6302.8
$ws.Range("I65").Value = 5737.875
$ws.Range("J65").Value = 6948.4287
$ws.Range("K65").Value = 28689.375
$ws.Range("L65").Value = 34742.14350000001
$ws.Range("M65").Value = -25569.375
$ws.Range("N65").Value = -40982.14350000001

# row 100 on ALC (-5625,22 +5628,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3997.3
$ws.Range("I100").Value = 2493.25
$ws.Range("K100").Value = 2493.25
$ws.Range("M100").Value = -1952.25

# row 132 on ALC (-7220,22 +7223,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1822.0465
$ws.Range("I132").Value = 1598.85
$ws.Range("K132").Value = 4796.549999999999
$ws.Range("M132").Value = -2266.549999999999

# row 138 on ALC (-7523,25 +7526,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2640.7878
$ws.Range("J138").Value = 2728.152
$ws.Range("L138").Value = 8184.456
$ws.Range("N138").Value = -18464.456

# row 32 on ARM (-9295,22 +9298,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19609938
$ws.Range("I32").Value = 19722714
$ws.Range("K32").Value = 19722714
$ws.Range("M32").Value = -19722427

# row 61 on ARM (-10707,25 +10710,25)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3456.5356
$ws.Range("I61").Value = 3246.389
$ws.Range("J61").Value = 3834.8
$ws.Range("K61").Value = 3246.389
$ws.Range("L61").Value = 3834.8
$ws.Range("M61").Value = -3034.389
$ws.Range("N61").Value = -4258.8

# row 63 on ARM (-10808,22 +10811,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3907.5557
$ws.Range("I63").Value = 2499
$ws.Range("K63").Value = 2499
$ws.Range("M63").Value = -1813

# row 66 on ARM (-10955,22 +10958,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3907.5557
$ws.Range("I66").Value = 2499
$ws.Range("K66").Value = 12495
$ws.Range("M66").Value = -9063

# row 74 on ARM (-11335,22 +11338,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2840.2173
$ws.Range("I74").Value = 2815.6
$ws.Range("K74").Value = 2815.6
$ws.Range("M74").Value = -1941.6

# row 77 on ARM (-11485,22 +11488,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2840.2173
$ws.Range("I77").Value = 2815.6
$ws.Range("K77").Value = 14078
$ws.Range("M77").Value = -9710

# row 97 on ARM (-12465,25 +12468,25)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2819.25
$ws.Range("I97").Value = 2525.6667
$ws.Range("J97").Value = 3700
$ws.Range("K97").Value = 2525.6667
$ws.Range("L97").Value = 3700
$ws.Range("M97").Value = -2029.6667
$ws.Range("N97").Value = -4692

# row 134 on ARM (-14284,22 +14287,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 119999.5
$ws.Range("J134").Value = 119999.5
$ws.Range("L134").Value = 119999.5
$ws.Range("N134").Value = -130139.5

# row 136 on ARM (-14385,25 +14388,25)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3456.5356
$ws.Range("I136").Value = 3246.389
$ws.Range("J136").Value = 3834.8
$ws.Range("K136").Value = 9739.167000000001
$ws.Range("L136").Value = 11504.4
$ws.Range("M136").Value = -7189.167000000001
$ws.Range("N136").Value = -16604.4

# row 105 on BSM (-19775,22 +19778,22)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2349.9333
$ws.Range("I105").Value = 1444.25
$ws.Range("K105").Value = 1444.25
$ws.Range("M105").Value = 302.75

# row 134 on BSM (-21190,22 +21193,22)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1744090.2
$ws.Range("I134").Value = 1932102.8
$ws.Range("K134").Value = 5796308.4
$ws.Range("M134").Value = -5793773.4

# row 138 on BSM (-21383,22 +21386,22)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 123333
$ws.Range("J138").Value = 123333
$ws.Range("L138").Value = 123333
$ws.Range("N138").Value = -133613

# row 31 on CRP (-23094,25 +23097,25)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6316.3076
$ws.Range("J31").Value = 6617.6665
$ws.Range("L31").Value = 6617.6665
$ws.Range("N31").Value = -7207.6665

# row 34 on CRP (-23244,25 +23247,25)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6316.3076
$ws.Range("J34").Value = 6617.6665
$ws.Range("L34").Value = 6617.6665
$ws.Range("N34").Value = -7021.6665

# row 86 on CRP (-25744,25 +25747,25)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 38353
$ws.Range("I86").Value = 30569.857
$ws.Range("J86").Value = 39786.74
$ws.Range("K86").Value = 30569.857
$ws.Range("L86").Value = 39786.74
$ws.Range("M86").Value = -29446.857
$ws.Range("N86").Value = -42032.74

# row 89 on CRP (-25894,25 +25897,25)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 38353
$ws.Range("I89").Value = 30569.857
$ws.Range("J89").Value = 39786.74
$ws.Range("K89").Value = 152849.285
$ws.Range("L89").Value = 198933.7
$ws.Range("M89").Value = -147233.285
$ws.Range("N89").Value = -210165.7

# row 131 on CRP (-27976,22 +27979,22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 40934.5
$ws.Range("J131").Value = 40934.5
$ws.Range("L131").Value = 40934.5
$ws.Range("N131").Value = -51014.5

# row 5 on CUL (-28762,25 +28765,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1495.55
$ws.Range("I5").Value = 1371.625
$ws.Range("J5").Value = 1578.1666
$ws.Range("K5").Value = 4114.875
$ws.Range("L5").Value = 4734.4998
$ws.Range("M5").Value = -4002.875
$ws.Range("N5").Value = -4958.4998

# row 34 on CUL (-30219,25 +30222,22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 134.57143
$ws.Range("I34").Value = 134.57143
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 403.71429
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -319.71429
$ws.Range("N34").ClearContents()

# row 39 on CUL (-30467,22 +30467,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2450.75
$ws.Range("I39").Value = 800
$ws.Range("J39").Value = 3001
$ws.Range("K39").Value = 2400
$ws.Range("L39").Value = 9003
$ws.Range("N39").Value = -9591
$ws.Range("M39").Value = -2106

# row 55 on CUL (-31260,25 +31263,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1635.6923
$ws.Range("J55").Value = 3503.5
$ws.Range("L55").Value = 10510.5
$ws.Range("N55").Value = -10864.5

# row 60 on CUL (-31514,22 +31517,22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1002
$ws.Range("J60").Value = 1002
$ws.Range("L60").Value = 3006
$ws.Range("N60").Value = -3508

# row 92 on CUL (-33082,25 +33085,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1700
$ws.Range("J92").Value = 1500
$ws.Range("L92").Value = 4500
$ws.Range("N92").Value = -6996

# row 135 on CUL (-35240,25 +35243,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1495.55
$ws.Range("I135").Value = 1371.625
$ws.Range("J135").Value = 1578.1666
$ws.Range("K135").Value = 12344.625
$ws.Range("L135").Value = 14203.4994
$ws.Range("M135").Value = -9809.625
$ws.Range("N135").Value = -19273.4994

# row 52 on GSM (-38130,25 +38133,25)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 45063.332
$ws.Range("I52").Value = 46000
$ws.Range("J52").Value = 44595
$ws.Range("K52").Value = 46000
$ws.Range("L52").Value = 44595
$ws.Range("M52").Value = -45741
$ws.Range("N52").Value = -45113

# row 54 on GSM (-38228,22 +38231,22)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 32698.334
$ws.Range("J54").Value = 32698.334
$ws.Range("L54").Value = 32698.334
$ws.Range("N54").Value = -33478.334

# row 92 on GSM (-40060,22 +40063,22)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 250251
$ws.Range("J92").Value = 250251
$ws.Range("L92").Value = 250251
$ws.Range("N92").Value = -253995

# row 132 on GSM (-42023,25 +42026,25)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3795.2856
$ws.Range("I132").Value = 3511.1667
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 10533.5001
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -8003.500100000001
$ws.Range("N132").Value = -21560

# row 7 on LTW (-42855,25 +42858,25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6350.4375
$ws.Range("I7").Value = 5516.615
$ws.Range("J7").Value = 9963.666999999999
$ws.Range("K7").Value = 5516.615
$ws.Range("L7").Value = 9963.666999999999
$ws.Range("M7").Value = -5404.615
$ws.Range("N7").Value = -10187.667

# row 40 on LTW (-44472,22 +44475,22)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4552.067
$ws.Range("I40").Value = 3930.3
$ws.Range("K40").Value = 3930.3
$ws.Range("M40").Value = -3794.3

# row 87 on LTW (-46751,25 +46754,25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 98999.5
$ws.Range("J87").Value = 122999
$ws.Range("L87").Value = 122999
$ws.Range("N87").Value = -125245

# row 90 on LTW (-46901,25 +46904,25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 98999.5
$ws.Range("J90").Value = 122999
$ws.Range("L90").Value = 368997
$ws.Range("N90").Value = -380229

# row 126 on LTW (-48656,25 +48659,25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6350.4375
$ws.Range("I126").Value = 5516.615
$ws.Range("J126").Value = 9963.666999999999
$ws.Range("K126").Value = 16549.845
$ws.Range("L126").Value = 29891.001
$ws.Range("M126").Value = -14079.845
$ws.Range("N126").Value = -34831.001

# row 128 on LTW (-48757,22 +48760,22)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 111000
$ws.Range("J128").Value = 111000
$ws.Range("L128").Value = 111000
$ws.Range("N128").Value = -120960

# row 132 on LTW (-48950,25 +48953,25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3727.8774
$ws.Range("J132").Value = 3949.3
$ws.Range("L132").Value = 11847.9
$ws.Range("N132").Value = -16907.9

